$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.104.01"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.833.89"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.03"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6279"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07472"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2924"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.07"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07724"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "1.836.04"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.974"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6683"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.70"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009305"
$ws.Range("E16").Value = "  -6.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.034"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "29.125.98"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "2.084.31"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.59"
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "222.96"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.127"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.12"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.503"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.90"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.504"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05625"
$ws.Range("E30").Value = "  +8.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.146"
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.073"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7507"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.851"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.136"
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.612"
$ws.Range("E37").Value = "  -3.68%  "
$ws.Range("D38").Value = "1.231.60"
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.753"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01783"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.567"
$ws.Range("E41").Value = "  +3.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8957"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.003"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.93"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").Value = "1.987.73"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("B46").Value = "XinFinNetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07912"
$ws.Range("E46").Value = "  +15.91%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.59"
$ws.Range("E47").Value = "  +2.22%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000123"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5094"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4056"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.038"
$ws.Range("E51").Value = "  +1.93%  "
